$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tl = $s.TimeLine
$ms = $tl.MainSequence
Write-Host "Before count:" $ms.Count
for ($i = 1; $i -le $ms.Count; $i++) {
    $eff = $ms.Item($i)
    Write-Host $i ":" $eff.Shape.Id $eff.Shape.Name "Trigger=" $eff.Timing.TriggerType "Paragraph=" $eff.Paragraph
}

$shp = $s.Shapes.Item(2)  # Rectangle 5, spid 21507

# delete old item 2 and 3 (pRg1, pRg2 with previous)
$ms.Item(3).Delete()
$ms.Item(2).Delete()

Write-Host "After delete count:" $ms.Count
for ($i = 1; $i -le $ms.Count; $i++) {
    $eff = $ms.Item($i)
    Write-Host $i ":" $eff.Shape.Id $eff.Shape.Name "Trigger=" $eff.Timing.TriggerType "Paragraph=" $eff.Paragraph
}

# Re-add pRg1 as click after position 1
$e1 = $ms.AddEffect($shp, 10, 0, 1)
$e1.Paragraph = 2
$e1.MoveTo(2)

$e2 = $ms.AddEffect($shp, 10, 0, 1)
$e2.Paragraph = 3
$e2.MoveTo(3)

Write-Host "Final count:" $ms.Count
for ($i = 1; $i -le $ms.Count; $i++) {
    $eff = $ms.Item($i)
    Write-Host $i ":" $eff.Shape.Id $eff.Shape.Name "Trigger=" $eff.Timing.TriggerType "Paragraph=" $eff.Paragraph
}
